$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (B, C, E) ---
$ws.Columns("B").ColumnWidth = 40.15
$ws.Columns("C").ColumnWidth = 21.65
$ws.Columns("E").ColumnWidth = 42.98

# --- Row 51 (new row) build up first so new shared strings are created in the
#     same order as the authored workbook (E51, F51, G51 dates before the
#     processo numbers in C50/C51) ---
$ws.Range("A50").Copy($ws.Range("A51"))
$ws.Range("A51").Value = 49

$ws.Range("B49").Copy($ws.Range("B51"))

$ws.Range("D50").Copy($ws.Range("D51"))

# E51 / F51 need to be literal text (not auto-converted to dates)
$scratch = $ws.Range("AZ500")
$scratch.Formula = "=""07/08/2026"""
$scratch.Copy()
$ws.Range("E51").PasteSpecial(-4163)
$scratch.Clear()

$scratch.Formula = "=""12/08/2026"""
$scratch.Copy()
$ws.Range("F51").PasteSpecial(-4163)
$scratch.Clear()

# G51 (28/08/2026) is safe to set directly - "28" can't be a month
$ws.Range("G51").Value = "28/08/2026"

$ws.Range("H50").Copy($ws.Range("H51"))

# --- Row 50: change B50/C50 ---
# B50 becomes the same text as B49 ('202400834844), losing its quote-prefix style
$ws.Range("B49").Copy($ws.Range("B50"))
# C50 becomes a new quote-prefixed number
$ws.Range("C50").Value = "'2022111012772"

# --- C51 (new quote-prefixed number, created after C50's) ---
$ws.Range("C51").Value = "'202211101277"

# --- New empty styled cells (O45, E55) matching F57's style ---
$ws.Range("F57").Copy($ws.Range("O45"))
$ws.Range("F57").Copy($ws.Range("E55"))

# --- Selection / view ---
$ws.Range("E55").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
